# Regenerate the "想去人数" (wishlist / "want to go" count) figures in
# column F across all four sheets, matching the freshly scraped data
# snapshot (gh-pages output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 311
$ws.Cells.Item(6, 6).Value = 265
$ws.Cells.Item(9, 6).Value = 7420
$ws.Cells.Item(16, 6).Value = 1878
$ws.Cells.Item(18, 6).Value = 1282
$ws.Cells.Item(24, 6).Value = 7
$ws.Cells.Item(27, 6).Value = 507
$ws.Cells.Item(29, 6).Value = 4768
$ws.Cells.Item(31, 6).Value = 3968
$ws.Cells.Item(32, 6).Value = 2215
$ws.Cells.Item(33, 6).Value = 184
$ws.Cells.Item(37, 6).Value = 54
$ws.Cells.Item(38, 6).Value = 57
$ws.Cells.Item(44, 6).Value = 210
$ws.Cells.Item(45, 6).Value = 849
$ws.Cells.Item(47, 6).Value = 39

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(21, 6).Value = 4
$ws.Cells.Item(32, 6).Value = 1466
$ws.Cells.Item(33, 6).Value = 1466

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(8, 6).Value = 2931
$ws.Cells.Item(10, 6).Value = 1163
$ws.Cells.Item(13, 6).Value = 1879
$ws.Cells.Item(14, 6).Value = 8295
$ws.Cells.Item(15, 6).Value = 571

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 311
$ws.Cells.Item(7, 6).Value = 2931
$ws.Cells.Item(8, 6).Value = 7420
$ws.Cells.Item(10, 6).Value = 1163
$ws.Cells.Item(18, 6).Value = 1878
$ws.Cells.Item(20, 6).Value = 1282
$ws.Cells.Item(36, 6).Value = 4768
$ws.Cells.Item(38, 6).Value = 3968
$ws.Cells.Item(39, 6).Value = 184
$ws.Cells.Item(42, 6).Value = 54
$ws.Cells.Item(43, 6).Value = 57
$ws.Cells.Item(50, 6).Value = 1466
